# Adding screenshots as user input arg.
# Appends a third benchmark iteration ("iteration 2") as row 4 on Sheet1,
# and extends the bar chart's 7 series to include the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Append the new data row (row 4) ------------------------------------
$ws.Range("A4").Value = "iteration 2"
$ws.Cells.Item(4, 2).Value = 1443
$ws.Cells.Item(4, 3).Value = 399
$ws.Cells.Item(4, 4).Value = 181
$ws.Cells.Item(4, 5).Value = 30
$ws.Cells.Item(4, 6).Value = 2716
$ws.Cells.Item(4, 7).Value = 4330
$ws.Cells.Item(4, 8).Value = 5139

# Row 2/3's label cell (column A) carries the bold/centered/bordered style;
# replicate it onto A4 so the new label matches the existing ones (style index 1).
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Extend each chart series to cover the new row (A2:A4 / X2:X4) ------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$cols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $cols.Count; $i++) {
    $col = $cols[$i]
    $series = $chart.SeriesCollection($i + 1)
    $series.Formula = "=SERIES(Sheet1!`$$col`$1,Sheet1!`$A`$2:`$A`$4,Sheet1!`$$col`$2:`$$col`$4,$($i + 1))"
}
